$d = $word.ActiveDocument

# Fix the duplicated "valutata / valutata" wording:
# "), questo Comando propone che l'istanza sia valutata con parere contrario ..."
# becomes
# "), questo Comando propone parere contrario ..."
$apos = [char]0x2019
$old = "che l" + $apos + "istanza sia valutata con "
$found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $found) {
    Write-Host "WARNING: target phrase not found for replacement"
}

# Rename the "Collegamento Internet" character style to its English built-in
# name "Hyperlink".
$s = $d.Styles("Collegamento Internet")
$s.NameLocal = "Hyperlink"
